# Correction in SA algorithm and 746 logs
# Updates the Fitness (column C) values for run_3.xlsx as produced by the
# corrected simulated-annealing logging. Values are set in contiguous
# row-range blocks that share the same corrected fitness value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C7").Value = 8132
$ws.Range("C8:C14").Value = 8056
$ws.Range("C15:C19").Value = 7735
$ws.Range("C20:C32").Value = 7610
$ws.Range("C33:C68").Value = 7312
$ws.Range("C69:C207").Value = 7310
